# Regenerate save_data to use K (strikeouts) instead of Strike# in column G,
# after recomputing std/mean and the associated s_vals upstream.
# Here we just (re)write the freshly computed K values for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 0
    27 = 2
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 2
    33 = 0
    34 = 1
    35 = 0
    36 = 2
    37 = 0
    38 = 2
    39 = 0
    40 = 0
    41 = 1
    42 = 2
    43 = 2
    44 = 0
    45 = 1
    46 = 0
    47 = 3
    48 = 4
    49 = 3
    50 = 1
    51 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
